# Introduce and test behavior of a date column ("Birthdate") on the Person sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the whole new column (header + data) as a short date up front, like
# the rest of the sheet's per-column formatting (e.g. SuperUser's boolean
# format), so the date values below don't pick up a transient default format.
$ws.Range("G1:G3").NumberFormat = "dd/mm/yy"

# Header for the new column G.
$ws.Range("G1").Value = "Birthdate"

# Birthdates for the two existing people (row 2: Hugo Boss, row 3: Fritz Lang).
$ws.Range("G2").Value = (Get-Date -Year 1990 -Month 11 -Day 16 -Hour 0 -Minute 0 -Second 0)
$ws.Range("G3").Value = (Get-Date -Year 1992 -Month 12 -Day 17 -Hour 0 -Minute 0 -Second 0)

# Leave the selection on the newly added cell, matching interactive entry.
$ws.Range("G3").Select() | Out-Null
